$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Rows that share the "Ready for handoff" generate/handoff timestamps and get
# the "ht" priority assigned during report generation for handoff.
$rows = @(7, 8, 10, 11, 12, 13)

foreach ($r in $rows) {
    # "Latest HO Xliff Generate Date" column on Overview sheet (column G)
    $wsOverview.Range("G$r").Value = "2016-08-28 04:21:53"

    # "Latest Handoff Datetime" column on zh-cn sheet (column H)
    $wsZhCn.Range("H$r").Value = "2016-08-28 04:21:48"

    # "Latest Handoff Datetime" column on de-de sheet (column H)
    $wsDeDe.Range("H$r").Value = "2016-08-28 04:21:53"

    # Priority column (E) set to "ht" on zh-cn and de-de sheets
    $wsZhCn.Range("E$r").Value = "ht"
    $wsDeDe.Range("E$r").Value = "ht"
}
